$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 and 3 (the two "TID2" entries with Record ID 0), shifting
# everything below up by two rows.
$ws.Range("A2:XFD3").EntireRow.Delete() | Out-Null

# Leave the same kind of selection Excel would show after a row-delete of
# that range: the two rows that used to occupy 2:3 (now holding what was
# previously rows 4:5).
$ws.Range("A2:XFD3").Select() | Out-Null
